# Applies the crypto-price / volume refresh described in the commit
# "Updated cryptos list on Tue Nov 14 18:43:50 UTC 2023 with GitHub Actions".
#
# Sheet1 columns: A=rank (unchanged) B=Coin C=Link D=Price E=Volume(1h).
# Every data cell in the sheet is stored as text, so all writes below are
# text writes too. For the D-column values that look like a plain number
# (e.g. "241.81"), a bare .Value assignment would make Excel coerce the cell
# to a number (silently dropping formatting such as the trailing zero in
# "0.0730"). We force text the same way a user typing into the grid would -
# a leading apostrophe - and then ClearFormats() to drop the transient
# "stored as text" quote-prefix marker Excel leaves on the cell style, so
# the cell ends up identical to a plain text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.918.89'
$ws.Range("E2").Value = '  -2.36%  '
$ws.Range("D3").Value = '2.007.70'
$ws.Range("E3").Value = '  -3.99%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '''241.81'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = '''0.647'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '''54.17'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").Value = '''58.21'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("D10").Value = '''0.357'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("D11").Value = '''0.0730'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.21%  '
$ws.Range("E12").Value = '  -4.38%  '
$ws.Range("D13").Value = '''0.881'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("D14").Value = '''13.99'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -7.21%  '
$ws.Range("D15").Value = '2.312.65'
$ws.Range("E15").Value = '  -3.45%  '
$ws.Range("D16").Value = '''5.22'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.43%  '
$ws.Range("D17").Value = '2.004.59'
$ws.Range("E17").Value = '  -4.68%  '
$ws.Range("D18").Value = '35.812.79'
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("D19").Value = '''16.95'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").Value = '''70.64'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.91%  '
$ws.Range("D21").Value = '0.0₃0841'
$ws.Range("E21").Value = '  -4.86%  '
$ws.Range("D22").Value = '''235.14'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("D23").Value = '''5.08'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -6.72%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").Value = '''2.33'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.58%  '
$ws.Range("D26").Value = '''2.23'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.86%  '
$ws.Range("D27").Value = '''9.10'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -6.62%  '
$ws.Range("D28").Value = '''162.50'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("D29").Value = '''19.57'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.20%  '
$ws.Range("D30").Value = '''0.119'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.21%  '
$ws.Range("E31").Value = '  -3.05%  '
$ws.Range("D32").Value = '''4.83'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -9.12%  '
$ws.Range("D33").Value = '''0.0589'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.48%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '''4.27'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -9.26%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '''0.0888'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +7.20%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  -1.27%  '
$ws.Range("D38").Value = '''2.15'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -11.01%  '
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("B40").Value = 'HuobiToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D40").Value = '''2.88'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''1.18'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -7.05%  '
$ws.Range("D42").Value = '''0.0211'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.22%  '
$ws.Range("E43").Value = '  -6.26%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = '''0.0886'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -7.11%  '
$ws.Range("D45").Value = '1.377.31'
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''90.37'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.19%  '
$ws.Range("D47").Value = '''7.31'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").Value = '''15.28'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.91%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = '''2.22'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -9.75%  '
$ws.Range("D51").Value = '''45.13'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.24%  '
